$d = $word.ActiveDocument

$replacements = @(
    @("11×63=693", "96×66=6336"),
    @("84×49=4116", "56×38=2128"),
    @("11×11=121", "44×70=3080"),
    @("32×99=3168", "76×42=3192"),
    @("21×42=882", "85×72=6120"),
    @("82×38=3116", "37×29=1073"),
    @("52×21=1092", "92×14=1288"),
    @("44×80=3520", "69×91=6279"),
    @("51×64=3264", "12×66=792"),
    @("70×35=2450", "77×24=1848"),
    @("21×94=1974", "12×61=732"),
    @("91×14=1274", "72×26=1872"),
    @("58×97=5626", "60×45=2700"),
    @("78×59=4602", "75×32=2400"),
    @("25×32=800", "73×71=5183"),
    @("18×26=468", "69×53=3657"),
    @("63×41=2583", "77×28=2156"),
    @("39×61=2379", "36×74=2664"),
    @("59×73=4307", "20×38=760"),
    @("85×46=3910", "70×52=3640"),
    @("36×64=2304", "20×13=260"),
    @("63×13=819", "96×74=7104"),
    @("99×12=1188", "99×24=2376"),
    @("65×51=3315", "47×19=893"),
    @("21×40=840", "70×74=5180")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
